# Fix occurrence field in standard import
# Replace `Encounter.occurrenceID` with `Encounter.sightingID` on the
# "Import Sheet" header row (cell C1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import Sheet")

$ws.Range("C1").Value = "Encounter.sightingID"
